# Auto-generated script to apply scheduled market-data refresh to Mandragora Profits workbook.
# Updates currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ / LevePriceNQ / LevePriceHQ /
# LeveProfitNQ / LeveProfitHQ (columns H-N) for the rows whose market data changed.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 15336.308
$ws.Range("I62").Value = 3997.625
$ws.Range("K62").Value = 3997.625
$ws.Range("M62").Value = -3373.625
$ws.Range("H65").Value = 15336.308
$ws.Range("I65").Value = 3997.625
$ws.Range("K65").Value = 19988.125
$ws.Range("M65").Value = -16868.125
$ws.Range("H106").Value = 2099
$ws.Range("I106").Value = 1748.75
$ws.Range("J106").Value = 3500
$ws.Range("K106").Value = 1748.75
$ws.Range("L106").Value = 3500
$ws.Range("M106").Value = -1117.75
$ws.Range("N106").Value = -4762
$ws.Range("H135").Value = 822.4545000000001
$ws.Range("I135").Value = 782.8889
$ws.Range("J135").Value = 1000.5
$ws.Range("K135").Value = 7046.0001
$ws.Range("L135").Value = 9004.5
$ws.Range("M135").Value = -4511.0001
$ws.Range("N135").Value = -14074.5
$ws.Range("H137").Value = 3075.647
$ws.Range("I137").Value = 4735.75
$ws.Range("J137").Value = 1600
$ws.Range("K137").Value = 14207.25
$ws.Range("L137").Value = 4800
$ws.Range("M137").Value = -11657.25
$ws.Range("N137").Value = -9900

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3353.0667
$ws.Range("I2").Value = 2191.3333
$ws.Range("K2").Value = 2191.3333
$ws.Range("M2").Value = -2078.3333
$ws.Range("H32").Value = 7654.477
$ws.Range("I32").Value = 8147.4653
$ws.Range("J32").Value = 3569.7144
$ws.Range("K32").Value = 8147.4653
$ws.Range("L32").Value = 3569.7144
$ws.Range("M32").Value = -7860.4653
$ws.Range("N32").Value = -4143.7144
$ws.Range("H102").Value = 2590
$ws.Range("I102").Value = 2531.4285
$ws.Range("J102").Value = 3000
$ws.Range("K102").Value = 2531.4285
$ws.Range("L102").Value = 3000
$ws.Range("M102").Value = -909.4285
$ws.Range("N102").Value = -6244
$ws.Range("H110").Value = 2441.1
$ws.Range("I110").Value = 1037.1818
$ws.Range("J110").Value = 4157
$ws.Range("K110").Value = 1037.1818
$ws.Range("L110").Value = 4157
$ws.Range("M110").Value = 1007.8182
$ws.Range("N110").Value = -8247
$ws.Range("H116").Value = 3353.0667
$ws.Range("I116").Value = 2191.3333
$ws.Range("K116").Value = 2191.3333
$ws.Range("M116").Value = 102.6667000000002
$ws.Range("H132").Value = 3940.4285
$ws.Range("I132").Value = 1551.6177
$ws.Range("J132").Value = 9355.066000000001
$ws.Range("K132").Value = 4654.8531
$ws.Range("L132").Value = 28065.198
$ws.Range("M132").Value = -2124.8531
$ws.Range("N132").Value = -33125.198

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3353.0667
$ws.Range("I3").Value = 2191.3333
$ws.Range("K3").Value = 2191.3333
$ws.Range("M3").Value = -2077.3333
$ws.Range("H99").Value = 2650.7026
$ws.Range("I99").Value = 2286.0417
$ws.Range("J99").Value = 3323.923
$ws.Range("K99").Value = 2286.0417
$ws.Range("L99").Value = 3323.923
$ws.Range("M99").Value = -788.0417000000002
$ws.Range("N99").Value = -6319.923
$ws.Range("H134").Value = 4083.9321
$ws.Range("I134").Value = 1813.1177
$ws.Range("J134").Value = 7172.24
$ws.Range("K134").Value = 5439.3531
$ws.Range("L134").Value = 21516.72
$ws.Range("M134").Value = -2904.3531
$ws.Range("N134").Value = -26586.72

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4496.091
$ws.Range("I16").Value = 3932.125
$ws.Range("K16").Value = 3932.125
$ws.Range("M16").Value = -3645.125
$ws.Range("H110").Value = 0
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("M110").ClearContents()
$ws.Range("N110").ClearContents()
$ws.Range("H113").Value = 4496.091
$ws.Range("I113").Value = 3932.125
$ws.Range("K113").Value = 3932.125
$ws.Range("M113").Value = -1762.125
$ws.Range("H132").Value = 3114.55
$ws.Range("I132").Value = 2213.2856
$ws.Range("J132").Value = 3599.8462
$ws.Range("K132").Value = 6639.8568
$ws.Range("L132").Value = 10799.5386
$ws.Range("M132").Value = -4109.8568
$ws.Range("N132").Value = -15859.5386
$ws.Range("H134").Value = 1691.3478
$ws.Range("I134").Value = 857.24243
$ws.Range("J134").Value = 3808.6924
$ws.Range("K134").Value = 2571.72729
$ws.Range("L134").Value = 11426.0772
$ws.Range("M134").Value = -36.72728999999981
$ws.Range("N134").Value = -16496.0772

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1388
$ws.Range("I68").Value = 1213.3334
$ws.Range("J68").Value = 1446.2222
$ws.Range("K68").Value = 3640.0002
$ws.Range("L68").Value = 4338.6666
$ws.Range("M68").Value = -2829.0002
$ws.Range("N68").Value = -5960.6666
$ws.Range("H71").Value = 1388
$ws.Range("I71").Value = 1213.3334
$ws.Range("J71").Value = 1446.2222
$ws.Range("K71").Value = 10920.0006
$ws.Range("L71").Value = 13015.9998
$ws.Range("M71").Value = -6864.000599999999
$ws.Range("N71").Value = -21127.9998
$ws.Range("H80").Value = 2684.2144
$ws.Range("I80").Value = 1049.75
$ws.Range("J80").Value = 3338
$ws.Range("K80").Value = 3149.25
$ws.Range("L80").Value = 10014
$ws.Range("M80").Value = -2213.25
$ws.Range("N80").Value = -11886
$ws.Range("H83").Value = 2684.2144
$ws.Range("I83").Value = 1049.75
$ws.Range("J83").Value = 3338
$ws.Range("K83").Value = 9447.75
$ws.Range("L83").Value = 30042
$ws.Range("M83").Value = -4767.75
$ws.Range("N83").Value = -39402
$ws.Range("H98").Value = 1137.5625
$ws.Range("I98").Value = 417
$ws.Range("J98").Value = 1569.9
$ws.Range("K98").Value = 1251
$ws.Range("L98").Value = 4709.700000000001
$ws.Range("M98").Value = 247
$ws.Range("N98").Value = -7705.700000000001
$ws.Range("H107").Value = 426
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 426
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 1278
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -5118
$ws.Range("H122").Value = 3386.82
$ws.Range("J122").Value = 3638.0652
$ws.Range("L122").Value = 32742.5868
$ws.Range("N122").Value = -37642.5868

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1737959.1
$ws.Range("I132").Value = 2084869.6
$ws.Range("K132").Value = 6254608.800000001
$ws.Range("M132").Value = -6252078.800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2247.8572
$ws.Range("I68").Value = 2086.6667
$ws.Range("J68").Value = 2538
$ws.Range("K68").Value = 2086.6667
$ws.Range("L68").Value = 2538
$ws.Range("M68").Value = -1337.6667
$ws.Range("N68").Value = -4036
$ws.Range("H71").Value = 2247.8572
$ws.Range("I71").Value = 2086.6667
$ws.Range("J71").Value = 2538
$ws.Range("K71").Value = 10433.3335
$ws.Range("L71").Value = 12690
$ws.Range("M71").Value = -6689.333500000001
$ws.Range("N71").Value = -20178
$ws.Range("H136").Value = 62501524
$ws.Range("I136").Value = 71430030
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 214290090
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -214287540
$ws.Range("N136").Value = -11100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 6251399.5
$ws.Range("I136").Value = 13889648
$ws.Range("J136").Value = 1923.7273
$ws.Range("K136").Value = 41668944
$ws.Range("L136").Value = 5771.1819
$ws.Range("N136").Value = -10871.1819
